# #5: cash & deposit done
#
# The 存款 (deposit) sheet's header row had accidentally been filled with
# a copy of the first data row's values instead of real column headers.
# Fix the header row (B1:F1) to proper labels, and extend the sheet with
# the bank/deposit_type/currency-style metadata columns (G:M) that the
# other property sheets (土地/建物/汽車) already carry: property_category,
# category, date, legislator_name, legislator_id, source_file, index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) ---
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# New header cells G1:M1 - give them the same bold/bordered look as the
# existing header cells before filling in their text.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows (2-5): new columns G:M only; B:F are untouched ---
for ($r = 2; $r -le 5; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    # Leading apostrophe forces this to stay plain text (matching the
    # other sheets, where dates like this are stored as text, not real
    # date serials) instead of Excel auto-parsing it as a date.
    $ws.Cells.Item($r, 9).Value = "'2012-04-26"
    $ws.Cells.Item($r, 10).Value = "楊瓊瓔"
    $ws.Cells.Item($r, 11).Value = 854
    $ws.Cells.Item($r, 12).Value = "tmp8a701"
    $ws.Cells.Item($r, 13).Value = $idx
}
